# Insert a new data row at row 458 (pushing existing rows 458-573 down to 459-574)
# and populate it with the new record's values, as described by the diff:
#   - dimension grows from A1:R573 to A1:R574
#   - a new "Acelga" price record is inserted right after row 457

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 458, shifting rows 458:573 down to 459:574.
$ws.Rows.Item(458).Insert()

# Populate the newly inserted row 458 with the new record's values.
$ws.Range("A458").Value = 10
$ws.Range("B458").Value = "Vega Modelo de Temuco"
$ws.Range("C458").Value = "La Araucanía"
$ws.Range("D458").Value = 45204
$ws.Range("E458").Value = 9
$ws.Range("F458").Value = 100112009
$ws.Range("G458").Value = "Acelga"
$ws.Range("H458").Value = "Sin especificar"
$ws.Range("I458").Value = "Primera"
$ws.Range("J458").Value = 55
$ws.Range("K458").Value = 7000
$ws.Range("L458").Value = 7000
$ws.Range("M458").Value = 7000
$ws.Range("N458").Value = "$/docena de atados (12 kilos)"
$ws.Range("O458").Value = "Región del Maule"
$ws.Range("P458").Value = 583
$ws.Range("Q458").Value = 12
$ws.Range("R458").Value = "Hortaliza"
